$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column G (ExpectedFilenames) with the new report-name patterns ---
# Existing rows 2-10 get new text; rows 11-13 are brand new.
$ws.Range("G2").Value  = "StandardExcelReport-Takeda - MM Maintenance-Clinical-2023_"
$ws.Range("G3").Value  = "ExcelReport-Takeda-MM Maintenance-Clinical-"
$ws.Range("G4").Value  = "WordReport-Takeda - MM Maintenance-Clinical-"
$ws.Range("G5").Value  = "StandardExcelReport-Takeda - MM Maintenance-Economic-2023_"
$ws.Range("G6").Value  = "ExcelReport-Takeda-MM Maintenance-Economic-"
$ws.Range("G7").Value  = "WordReport-Takeda - MM Maintenance-Economic-"
$ws.Range("G8").Value  = "StandardExcelReport-Takeda - MM Maintenance-Quality of Life-2023_"
$ws.Range("G9").Value  = "ExcelReport-Takeda-MM Maintenance-Quality of Life-"
$ws.Range("G10").Value = "WordReport-Takeda - MM Maintenance-Quality of Life-"
$ws.Range("G11").Value = "StandardExcelReport-Takeda - MM Maintenance-Real-world Evidence-2023_"
$ws.Range("G12").Value = "ExcelReport-Takeda-MM Maintenance-Real-world Evidence-"
$ws.Range("G13").Value = "WordReport-Takeda - MM Maintenance-Real-world Evidence-"

# --- Strip the stray no-op "applyAlignment" style off the data cells that had it ---
$ws.Range("A2:A5").Style = "Normal"
$ws.Range("D3:D5").Style = "Normal"
$ws.Range("E2:F5").Style = "Normal"
$ws.Range("G2:G13").Style = "Normal"

# --- Update the view: scrolled to column F, selection now G2:G13 ---
$ws.Range("G2:G13").Select()
$excel.ActiveWindow.ScrollColumn = 6
